# Apply "Gave default value to the pads" edit:
# Each of the four `range(<x>._<y>Pad)` template expressions gets an
# `or 80` default value inserted right before the closing parenthesis,
# i.e. `range(work._startDatePad)` -> `range(work._startDatePad or 80)`.

$d = $word.ActiveDocument

# 1) Insert the extra (structurally empty) paragraph that appears right
#    before the "{% for work in res.work -%}" paragraph in the WORK
#    EXPERIENCE section.
$marker = $d.Content
$found = $marker.Find.Execute("{% for work in res.work -%}", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "", 0)
if ($found) {
    $insertionPoint = $d.Range($marker.Start, $marker.Start)
    $insertionPoint.InsertParagraphBefore()
}

# 2) Give each of the four "pad" variables a default value of 80 by
#    inserting " or 80" right before the closing parenthesis of the
#    enclosing range(...) call.
$pads = @("work._startDatePad", "work._locationPad", "education._endDatePad", "education._locationPad")
foreach ($pad in $pads) {
    $searchText = $pad + ")"
    $rng = $d.Content
    $ok = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, `
                             $true, 1, $false, "", 0)
    if ($ok) {
        $insertAt = $d.Range($rng.End - 1, $rng.End - 1)
        $insertAt.InsertBefore(" or 80")
    }
}
